$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.723.08"
$ws.Range("E2").Value = "  +2.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.572.77"
$ws.Range("E3").Value = "  +2.03%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.41"
$ws.Range("E5").Value = "  +2.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.41"
$ws.Range("E6").Value = "  +2.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.566.63"
$ws.Range("E7").Value = "  +2.05%  "

$ws.Range("E8").Value = "  +1.41%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("E10").Value = "  +6.52%  "

$ws.Range("E11").Value = "  +8.59%  "

$ws.Range("E12").Value = "  +2.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.72"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("E14").Value = "  +1.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.145.23"
$ws.Range("E15").Value = "  +1.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.41"
$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "612.78"
$ws.Range("E17").Value = "  +0.84%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.582.03"
$ws.Range("E18").Value = "  +2.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.766.09"
$ws.Range("E19").Value = "  +2.88%  "

$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.45"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.886"
$ws.Range("E22").Value = "  +0.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.20"
$ws.Range("E23").Value = "  -17.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.84"
$ws.Range("E24").Value = "  +1.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.15"
$ws.Range("E25").Value = "  +1.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.78"
$ws.Range("E26").Value = "  -0.93%  "

$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("E28").Value = "  +1.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.00"
$ws.Range("E29").Value = "  +4.99%  "

$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("E31").Value = "  -0.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.07"
$ws.Range("E32").Value = "  -1.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.20"
$ws.Range("E33").Value = "  +4.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "663.07"
$ws.Range("E34").Value = "  +7.24%  "

$ws.Range("E35").Value = "  -0.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.71"
$ws.Range("E36").Value = "  +8.42%  "

$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.84"
$ws.Range("E38").Value = "  +1.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0479"
$ws.Range("E39").Value = "  +8.86%  "

$ws.Range("E40").Value = "  +0.53%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("E42").Value = "  +6.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.394.51"
$ws.Range("E43").Value = "  +0.89%  "

$ws.Range("E44").Value = "  -0.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0716"
$ws.Range("E45").Value = "  +3.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "32.97"
$ws.Range("E46").Value = "  +1.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.96"
$ws.Range("E47").Value = "  +8.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.67"
$ws.Range("E48").Value = "  +6.40%  "

$ws.Range("E49").Value = "  +1.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.24"
$ws.Range("E50").Value = "  -0.40%  "
